$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The BOM rows for the male pin headers (J3***, J4***, PWR) previously had
# no "source" link in column F like every other single/multi row BOM entry.
# Add the missing source URL (as a new shared string) to F23:F25, apply the
# same "Hyperlink" cell style used elsewhere in column F, and wire up a
# single multi-cell hyperlink (like the existing F2:F10 one) pointing at the
# TME pin header listing.

$url = "https://www.tme.eu/hr/en/details/zl211-40kg-s/pin-headers/connfly/ds1022-1-40ruf11/"

$ws.Range("F23:F25").Style = "Hyperlink"

$ws.Range("F23").Value = $url
$ws.Range("F24").Value = $url
$ws.Range("F25").Value = $url

$linkRange = $ws.Range("F23:F25")
$null = $ws.Hyperlinks.Add($linkRange, $url, "", "", $url)

# Move the active selection from where it was left (B35) to D30, matching
# the saved view state of the edited workbook.
$null = $ws.Range("D30").Select()
